# Update usernames on the "user registration" report sheets (Hoja1-Hoja6)
$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("Hoja1")
$ws1.Range("C2").Value = "jonreyusr101"
$ws1.Range("C3").Value = "jonescusr102"
$ws1.Range("C4").Value = "abireyusr103"

$ws2 = $wb.Worksheets.Item("Hoja2")
$ws2.Range("C2").Value = "jonreyusr104"
$ws2.Range("C3").Value = "jonescusr105"

$ws3 = $wb.Worksheets.Item("Hoja3")
$ws3.Range("C2").Value = "jonreyusr106"
$ws3.Range("C3").Value = "jonescusr107"

$ws4 = $wb.Worksheets.Item("Hoja4")
$ws4.Range("C2").Value = "jonreyusr108"
$ws4.Range("C3").Value = "jonescusr109"
$ws4.Range("C4").Value = "abireyusr110"

$ws5 = $wb.Worksheets.Item("Hoja5")
$ws5.Range("C2").Value = "jonreyusr111"
$ws5.Range("C3").Value = "jonescusr112"

$ws6 = $wb.Worksheets.Item("Hoja6")
$ws6.Range("C2").Value = "jonreyusr113"
$ws6.Range("C3").Value = "jonescusr114"

$ws7 = $wb.Worksheets.Item("Hoja7")

# Update cell selections (cursor position) on each sheet as left by the author
$ws1.Range("C4").Select()
$ws3.Range("C3").Select()
$ws4.Range("C4").Select()
$ws5.Range("C2").Select()
$ws6.Range("C3").Select()
$ws7.Range("A3").Select()

# Hoja2 becomes the active tab (was Hoja16); set its selection too
$ws2.Activate()
$ws2.Range("C2").Select()
